$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.446.90"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -4.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.373.04"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -5.20%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "500.57"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.28"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.551"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.394.04"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0957"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.64%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.320"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.64"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -10.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.797.12"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -5.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.940.31"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.52"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.71%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.389.58"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.17"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "312.12"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.03"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.18"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.57"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.43%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.370"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -9.56%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -6.14%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "173.96"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.84%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.66"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.33%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0711"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.10"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.84%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -6.68%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.74"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.22"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.74"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.86"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.42"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -6.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.774"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -6.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "129.54"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.35"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.80"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.571"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "253.52"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -7.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0897"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0485"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.82"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0206"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -5.09%  "
